# Weekly data refresh: insert this week's "Frutilla" price rows at the top
# of the Mercado Mayorista Lo Valledor de Santiago block (before the oldest
# previously-known rows), pushing the rest of the historical rows down.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new blank rows at 525:527 - everything from row 525 downward
# (through the old last row 607) shifts down to 528:610.
$ws.Rows("525:527").Insert()

# Columns that are constant across every row of this block.
$ws.Range("A525:A527").Value = 6
$ws.Range("B525:B527").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C525:C527").Value = "Metropolitana"
$ws.Range("D525:D527").Value = 44474
$ws.Range("E525:E527").Value = 13
$ws.Range("F525:F527").Value = "Fruta"
$ws.Range("G525:G527").Value = 100101
$ws.Range("H525:H527").Value = "Berries"
$ws.Range("I525:I527").Value = 100112025
$ws.Range("J525:J527").Value = "Frutilla"
$ws.Range("K525:K527").Value = "Sin especificar"
$ws.Range("Q525:Q527").Value = "$/bandeja 7 kilos"
$ws.Range("R525:R527").Value = "Provincia de Melipilla"
$ws.Range("T525:T527").Value = 7

# Row 525: Especial
$ws.Range("L525").Value = "Especial"
$ws.Range("M525").Value = 800
$ws.Range("N525").Value = 12000
$ws.Range("O525").Value = 14000
$ws.Range("P525").Value = 13000
$ws.Range("S525").Value = 1857

# Row 526: Primera
$ws.Range("L526").Value = "Primera"
$ws.Range("M526").Value = 1000
$ws.Range("N526").Value = 9000
$ws.Range("O526").Value = 10000
$ws.Range("P526").Value = 9500
$ws.Range("S526").Value = 1357

# Row 527: Segunda
$ws.Range("L527").Value = "Segunda"
$ws.Range("M527").Value = 500
$ws.Range("N527").Value = 5000
$ws.Range("O527").Value = 6000
$ws.Range("P527").Value = 5500
$ws.Range("S527").Value = 786
